$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Units")

# Rename the "Cost_OM" / "Cost_var" table header cells on the Units sheet.
# These cells are the header row of the Excel Table ("Table1"), so editing
# the header cell value also renames the corresponding table column.
$ws.Range("S1").Value = "Cost_FOM (yearly)"
$ws.Range("T1").Value = "Cost_VOM (per unit)"

# New data entered under the renamed columns.
$ws.Range("S2").Value = 100
$ws.Range("T3").Value = 1

# Widen the two renamed columns so the new, longer headers are fully visible
# (mirrors Excel's automatic "best fit" column resize after a table header
# edit).
$ws.Columns.Item(19).AutoFit() | Out-Null
$ws.Columns.Item(20).AutoFit() | Out-Null

# Move the active selection, matching where the editor ended up after
# making the change.
$ws.Range("T4").Select() | Out-Null
